$wb = $excel.ActiveWorkbook

# --- 1. Update the "Generated from ..." timestamp on the Meta sheet ---
$meta = $wb.Worksheets.Item("Meta")
$meta.Range("B2").Value2 = "Generated from schema\openc2.jaen, Wed Mar 15 11:05:41 2017"

# --- 2. Add a new "actuators" entry to the "commands" vocabulary table on the Vocab sheet ---
# Before the edit, the "commands" vocabulary table (rows 104-109) looks like:
#   row 104: Vocabulary: commands
#   row 105: Description: Target used to query Actuator for its supported capabilities
#   row 107: Id | Value | Description
#   row 108: 1  | actions | results = JSON array of supported action verbs
#   row 109: 2  | schema  | results = JAEN syntax of supported commands
# We insert a new row 110 with: 3 | actuators | results = JSON array of actuator group names
$vocab = $wb.Worksheets.Item("Vocab")

# Insert a fresh row at position 110, pushing everything below down by one row.
$vocab.Rows.Item(110).Insert()

# Copy the formatting of the row above (the "schema" row) onto the new row so the
# new cells pick up the same borders / wrap / vertical alignment style.
$vocab.Range("B109:D109").Copy()
$vocab.Range("B110:D110").PasteSpecial(-4122)

# Fill in the new row's content.
$vocab.Range("B110").Value2 = 3
$vocab.Range("C110").Value2 = "actuators"
$vocab.Range("D110").Value2 = "results = JSON array of actuator group names"
